$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U1").Value = "Image GA Percobaan 1"
$ws.Range("V1").Value = "Image GA Percobaan 2"
$ws.Range("W1").Value = "Image GA Percobaan 3"
$ws.Range("X1").Value = "Image GA-ACO Percobaan 1"
$ws.Range("Y1").Value = "Image GA-ACO Percobaan 2"
$ws.Range("Z1").Value = "Image GA-ACO Percobaan 3"
$ws.Range("S2").Value = 0.0026
$ws.Range("T2").Value = 0.0007
$ws.Range("V2").Value = "./imageResult/t5_2_GA_10.png"
$ws.Range("W2").Value = "./imageResult/t5_3_GA_10.png"
$ws.Range("X2").Value = "./imageResult/t5_1_ACO_10.png"
$ws.Range("Y2").Value = "./imageResult/t5_2_ACO_10.png"
$ws.Range("Z2").Value = "./imageResult/t5_3_ACO_10.png"
$ws.Range("P3").Value = 0.001
$ws.Range("V3").Value = "./imageResult/t5_2_GA_50.png"
$ws.Range("W3").Value = "./imageResult/t5_3_GA_50.png"
$ws.Range("X3").Value = "./imageResult/t5_1_ACO_50.png"
$ws.Range("Y3").Value = "./imageResult/t5_2_ACO_50.png"
$ws.Range("Z3").Value = "./imageResult/t5_3_ACO_50.png"
$ws.Range("S4").Value = 0.0011
$ws.Range("V4").Value = "./imageResult/t5_2_GA_100.png"
$ws.Range("W4").Value = "./imageResult/t5_3_GA_100.png"
$ws.Range("X4").Value = "./imageResult/t5_1_ACO_100.png"
$ws.Range("Y4").Value = "./imageResult/t5_2_ACO_100.png"
$ws.Range("Z4").Value = "./imageResult/t5_3_ACO_100.png"
$ws.Range("I5").Value = 50.76352496678864
$ws.Range("J5").Value = 45.4099348702336
$ws.Range("K5").Value = 51.55699566478027
$ws.Range("L5").Value = 31.22691510942754
$ws.Range("M5").Value = 31.88252949105588
$ws.Range("V5").Value = "./imageResult/burma14_2_GA_10.png"
$ws.Range("W5").Value = "./imageResult/burma14_3_GA_10.png"
$ws.Range("X5").Value = "./imageResult/burma14_1_ACO_10.png"
$ws.Range("Y5").Value = "./imageResult/burma14_2_ACO_10.png"
$ws.Range("Z5").Value = "./imageResult/burma14_3_ACO_10.png"
$ws.Range("I6").Value = 37.15106686359928
$ws.Range("J6").Value = 44.8937944990004
$ws.Range("K6").Value = 39.6899265990123
$ws.Range("M6").Value = 31.88252949105588
$ws.Range("P6").Value = 0.0022
$ws.Range("Q6").Value = 0.0023
$ws.Range("S6").Value = 0.0018
$ws.Range("T6").Value = 0.0017
$ws.Range("V6").Value = "./imageResult/burma14_2_GA_50.png"
$ws.Range("W6").Value = "./imageResult/burma14_3_GA_50.png"
$ws.Range("X6").Value = "./imageResult/burma14_1_ACO_50.png"
$ws.Range("Y6").Value = "./imageResult/burma14_2_ACO_50.png"
$ws.Range("Z6").Value = "./imageResult/burma14_3_ACO_50.png"
$ws.Range("I7").Value = 34.83039736354373
$ws.Range("J7").Value = 43.35471087658466
$ws.Range("K7").Value = 36.56324197122565
$ws.Range("L7").Value = 32.11184821924887
$ws.Range("M7").Value = 31.88252949105588
$ws.Range("N7").Value = 31.22691510942754
$ws.Range("P7").Value = 0.0035
$ws.Range("Q7").Value = 0.0035
$ws.Range("R7").Value = 0.0025
$ws.Range("S7").Value = 0.0025
$ws.Range("T7").Value = 0.0024
$ws.Range("V7").Value = "./imageResult/burma14_2_GA_100.png"
$ws.Range("W7").Value = "./imageResult/burma14_3_GA_100.png"
$ws.Range("X7").Value = "./imageResult/burma14_1_ACO_100.png"
$ws.Range("Y7").Value = "./imageResult/burma14_2_ACO_100.png"
$ws.Range("Z7").Value = "./imageResult/burma14_3_ACO_100.png"
$ws.Range("I8").Value = 551695.2621313389
$ws.Range("J8").Value = 563773.6720051733
$ws.Range("K8").Value = 560316.2805824048
$ws.Range("L8").Value = 49143.7729793856
$ws.Range("N8").Value = 49215.61251916289
$ws.Range("O8").Value = 0.0091
$ws.Range("P8").Value = 0.0083
$ws.Range("Q8").Value = 0.008200000000000001
$ws.Range("R8").Value = 0.068
$ws.Range("S8").Value = 0.06569999999999999
$ws.Range("T8").Value = 0.06610000000000001
$ws.Range("V8").Value = "./imageResult/lin318_2_GA_10.png"
$ws.Range("W8").Value = "./imageResult/lin318_3_GA_10.png"
$ws.Range("X8").Value = "./imageResult/lin318_1_ACO_10.png"
$ws.Range("Y8").Value = "./imageResult/lin318_2_ACO_10.png"
$ws.Range("Z8").Value = "./imageResult/lin318_3_ACO_10.png"
$ws.Range("I9").Value = 533921.932031111
$ws.Range("J9").Value = 541321.7020989901
$ws.Range("K9").Value = 533837.8675633604
$ws.Range("L9").Value = 48382.86890044977
$ws.Range("M9").Value = 48382.86890044977
$ws.Range("N9").Value = 48382.86890044977
$ws.Range("O9").Value = 0.0335
$ws.Range("P9").Value = 0.0337
$ws.Range("Q9").Value = 0.0342
$ws.Range("R9").Value = 0.3146
$ws.Range("S9").Value = 0.316
$ws.Range("T9").Value = 0.3163
$ws.Range("V9").Value = "./imageResult/lin318_2_GA_50.png"
$ws.Range("W9").Value = "./imageResult/lin318_3_GA_50.png"
$ws.Range("X9").Value = "./imageResult/lin318_1_ACO_50.png"
$ws.Range("Y9").Value = "./imageResult/lin318_2_ACO_50.png"
$ws.Range("Z9").Value = "./imageResult/lin318_3_ACO_50.png"
$ws.Range("I10").Value = 487818.0359691684
$ws.Range("J10").Value = 502701.7910008261
$ws.Range("K10").Value = 523701.944058846
$ws.Range("L10").Value = 48108.8789343028
$ws.Range("M10").Value = 48783.31682921913
$ws.Range("N10").Value = 49215.61251916289
$ws.Range("O10").Value = 0.0663
$ws.Range("P10").Value = 0.0663
$ws.Range("Q10").Value = 0.0675
$ws.Range("R10").Value = 0.6301
$ws.Range("S10").Value = 0.6429
$ws.Range("T10").Value = 0.6337
$ws.Range("V10").Value = "./imageResult/lin318_2_GA_100.png"
$ws.Range("W10").Value = "./imageResult/lin318_3_GA_100.png"
$ws.Range("X10").Value = "./imageResult/lin318_1_ACO_100.png"
$ws.Range("Y10").Value = "./imageResult/lin318_2_ACO_100.png"
$ws.Range("Z10").Value = "./imageResult/lin318_3_ACO_100.png"
